# Update the crypto price/volume/name table to the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.152.80'
$ws.Range("E2").Value = '  -2.24%  '

$ws.Range("D3").Value = '1.575.38'
$ws.Range("E3").Value = '  -1.73%  '

$ws.Range("E4").Value = '  -0.46%  '

$ws.Range("D5").Value = "'" + '208.86'
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("E6").Value = '  -3.11%  '

$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  -1.65%  '

$ws.Range("E9").Value = '  -1.26%  '

$ws.Range("D10").Value = "'" + '19.55'
$ws.Range("E10").Value = '  -0.53%  '

$ws.Range("D11").Value = "'" + '0.0844'
$ws.Range("E11").Value = '  -0.38%  '

$ws.Range("D12").Value = '1.796.11'

$ws.Range("D13").Value = '1.594.94'
$ws.Range("E13").Value = '  -1.17%  '

$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").Value = "'" + '0.515'
$ws.Range("E15").Value = '  -1.98%  '

$ws.Range("D16").Value = "'" + '64.37'
$ws.Range("E16").Value = '  -1.04%  '

$ws.Range("D17").Value = '26.151.90'
$ws.Range("E17").Value = '  -2.16%  '

$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -1.99%  '

$ws.Range("D19").Value = "'" + '7.27'
$ws.Range("E19").Value = '  +1.73%  '

$ws.Range("D20").Value = "'" + '207.92'
$ws.Range("E20").Value = '  -0.83%  '

$ws.Range("E21").Value = '  -0.42%  '

$ws.Range("D22").Value = "'" + '4.26'
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("D23").Value = "'" + '2.18'
$ws.Range("E23").Value = '  -2.56%  '

$ws.Range("D24").Value = "'" + '8.83'
$ws.Range("E24").Value = '  -2.42%  '

$ws.Range("D25").Value = "'" + '143.72'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("E28").Value = '  -1.48%  '

$ws.Range("D29").Value = "'" + '15.20'
$ws.Range("E29").Value = '  -1.21%  '

$ws.Range("D30").Value = "'" + '0.0506'
$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("E32").Value = '  -2.10%  '

$ws.Range("E33").Value = '  +0.98%  '

$ws.Range("D34").Value = '1.278.49'
$ws.Range("E34").Value = '  -0.85%  '

$ws.Range("D35").Value = "'" + '0.614'
$ws.Range("E35").Value = '  +4.04%  '

$ws.Range("E36").Value = '  -1.50%  '

$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("D38").Value = "'" + '1.11'
$ws.Range("E38").Value = '  -10.16%  '

$ws.Range("E39").Value = '  -2.48%  '

$ws.Range("E40").Value = '  -2.16%  '

$ws.Range("E41").Value = '  -0.42%  '

$ws.Range("D42").Value = "'" + '5.58'
$ws.Range("E42").Value = '  +2.41%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = "'" + '0.763'
$ws.Range("E43").Value = '  -1.94%  '

$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = "'" + '2.13'
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("D45").Value = "'" + '62.37'

$ws.Range("D46").Value = '1.709.09'

$ws.Range("D47").Value = "'" + '88.75'
$ws.Range("E47").Value = '  -1.98%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  +0.72%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'" + '1.50'
$ws.Range("E49").Value = '  -3.70%  '

$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("E51").Value = '  -1.64%  '
